$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(154).Insert()

$ws.Range("A154").Value = 5
$ws.Range("B154").Value = "Macroferia Regional de Talca"
$ws.Range("C154").Value = "Maule"
$ws.Range("D154").Value = 44651
$ws.Range("E154").Value = 7
$ws.Range("F154").Value = 100112003
$ws.Range("G154").Value = "Ajo"
$ws.Range("H154").Value = "Chino"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 160
$ws.Range("K154").Value = 18000
$ws.Range("L154").Value = 18000
$ws.Range("M154").Value = 18000
$ws.Range("N154").Value = "`$/caja 10 kilos"
$ws.Range("O154").Value = "China"
$ws.Range("P154").Value = 1800
$ws.Range("Q154").Value = 10
$ws.Range("R154").Value = "Hortaliza"
